$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- AB column additions near header rows 12-14 ---
$ws.Range("AB12").Value = 20
$ws.Range("AB13").Value = 12
$ws.Range("AB14").Value = 12

# --- Z/AA/AB/AC table (grid size 15 by 15), rows 19-28 ---
$ws.Range("AA19").Value = 3
$ws.Range("AB19").Value = 11
$ws.Range("AC19").Value = 0

$ws.Range("AA20").Value = 5
$ws.Range("AB20").Value = 22
$ws.Range("AC20").Value = 1

$ws.Range("AA21").Value = 65
$ws.Range("AB21").Value = 52
$ws.Range("AC21").Value = 3

$ws.Range("AA22").Value = 2
$ws.Range("AB22").Value = 1
$ws.Range("AC22").Value = 0

$ws.Range("AA23").Value = 5
$ws.Range("AB23").Value = 2
$ws.Range("AC23").Value = 1

$ws.Range("AA24").Value = 28
$ws.Range("AB24").Value = 24
$ws.Range("AC24").Value = 1

$ws.Range("AA25").Value = 16
$ws.Range("AB25").Value = 8
$ws.Range("AC25").Value = 1

$ws.Range("AA26").Value = 5
$ws.Range("AB26").Value = 32
$ws.Range("AC26").Value = 1

$ws.Range("AA27").Value = 5
$ws.Range("AB27").Value = 2
$ws.Range("AC27").Value = 1

$ws.Range("AA28").Value = 9
$ws.Range("AB28").Value = 4
$ws.Range("AC28").Value = 1

# --- U/V/W/X table (grid size 5 by 5), first block rows 32-40 ---
$ws.Range("V32").Value = 1
$ws.Range("W32").Value = 10
$ws.Range("X32").Value = 1

$ws.Range("V33").Value = 2
$ws.Range("W33").Value = 21
$ws.Range("X33").Value = 1

$ws.Range("V34").Value = 1
$ws.Range("W34").Value = 0
$ws.Range("X34").Value = 1

$ws.Range("V35").Value = 3
$ws.Range("W35").Value = 11
$ws.Range("X35").Value = 1

$ws.Range("V36").Value = 1
$ws.Range("W36").Value = 0
$ws.Range("X36").Value = 1

$ws.Range("V37").Value = 1
$ws.Range("W37").Value = 0
$ws.Range("X37").Value = 1

$ws.Range("V38").Value = 1
$ws.Range("W38").Value = 0
$ws.Range("X38").Value = 1

$ws.Range("V39").Value = 2
$ws.Range("W39").Value = 1
$ws.Range("X39").Value = 1

$ws.Range("V40").Value = 4
$ws.Range("W40").Value = 22
$ws.Range("X40").Value = 1

# --- U/V/W/X table (grid size 5 by 5), second block rows 43-52 ---
$ws.Range("V43").Value = 7
$ws.Range("W43").Value = 3
$ws.Range("X43").Value = 0

$ws.Range("V44").Value = 15
$ws.Range("W44").Value = 17
$ws.Range("X44").Value = 1

$ws.Range("V45").Value = 32
$ws.Range("W45").Value = 16
$ws.Range("X45").Value = 1

$ws.Range("V46").Value = 6
$ws.Range("W46").Value = 3
$ws.Range("X46").Value = 1

$ws.Range("V47").Value = 11
$ws.Range("W47").Value = 5
$ws.Range("X47").Value = 2

$ws.Range("V48").Value = 19
$ws.Range("W48").Value = 9
$ws.Range("X48").Value = 1

$ws.Range("V49").Value = 26
$ws.Range("W49").Value = 13
$ws.Range("X49").Value = 1

$ws.Range("V50").Value = 163
$ws.Range("W50").Value = 91
$ws.Range("X50").Value = 2

$ws.Range("V51").Value = 112
$ws.Range("W51").Value = 56
$ws.Range("X51").Value = 1

$ws.Range("V52").Value = 11
$ws.Range("W52").Value = 5
$ws.Range("X52").Value = 2

# --- Selection / view state ---
$ws.Application.ActiveWindow.Zoom = 70
$ws.Range("AC28").Select()
